$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new reading was recorded for 2026/02/02 08:00 (月, hour 19 -> rank 194).
# Insert a fresh row right before the current row 770 (2026/12/29 entry) so
# every subsequent row shifts down by one, matching the sheet's chronological
# ordering (and growing the table from 811 to 812 data rows).
$ws.Rows.Item(770).Insert()

# Date column is stored as literal text (e.g. "2026/12/29"), not a real
# Excel date, in this sheet - force text entry with a leading apostrophe so
# it doesn't get auto-converted to a date serial number, then drop back to
# the sheet's default "Normal" style so no stray number-format is left on
# the cell.
$ws.Range("A770").Value = "'2026/02/02"
$ws.Range("A770").Style = "Normal"

$ws.Range("B770").Value = "月"
$ws.Range("C770").Value = 19
$ws.Range("D770").Value = 194
